$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("PIR")
$ws.Range("A75:A87").NumberFormat = "@"

$ws.Cells.Item(75, 1).Value = "2026-01-30"
$ws.Cells.Item(75, 2).Value = "16:55:02"
$ws.Cells.Item(75, 3).Value = "16:00"
$ws.Cells.Item(75, 4).Value = "Bathroom"
$ws.Cells.Item(75, 5).Value = "No Motion"
$ws.Cells.Item(75, 6).Value = "Inactive"

$ws.Cells.Item(76, 1).Value = "2026-01-30"
$ws.Cells.Item(76, 2).Value = "16:55:03"
$ws.Cells.Item(76, 3).Value = "16:00"
$ws.Cells.Item(76, 4).Value = "Bathroom"
$ws.Cells.Item(76, 5).Value = "No Motion"
$ws.Cells.Item(76, 6).Value = "Inactive"

$ws.Cells.Item(77, 1).Value = "2026-01-30"
$ws.Cells.Item(77, 2).Value = "16:55:07"
$ws.Cells.Item(77, 3).Value = "16:00"
$ws.Cells.Item(77, 4).Value = "Bathroom"
$ws.Cells.Item(77, 5).Value = "No Motion"
$ws.Cells.Item(77, 6).Value = "Inactive"

$ws.Cells.Item(78, 1).Value = "2026-01-30"
$ws.Cells.Item(78, 2).Value = "16:55:12"
$ws.Cells.Item(78, 3).Value = "16:00"
$ws.Cells.Item(78, 4).Value = "Bathroom"
$ws.Cells.Item(78, 5).Value = "No Motion"
$ws.Cells.Item(78, 6).Value = "Inactive"

$ws.Cells.Item(79, 1).Value = "2026-01-30"
$ws.Cells.Item(79, 2).Value = "16:55:17"
$ws.Cells.Item(79, 3).Value = "16:00"
$ws.Cells.Item(79, 4).Value = "Bathroom"
$ws.Cells.Item(79, 5).Value = "No Motion"
$ws.Cells.Item(79, 6).Value = "Inactive"

$ws.Cells.Item(80, 1).Value = "2026-01-30"
$ws.Cells.Item(80, 2).Value = "16:55:22"
$ws.Cells.Item(80, 3).Value = "16:00"
$ws.Cells.Item(80, 4).Value = "Bathroom"
$ws.Cells.Item(80, 5).Value = "No Motion"
$ws.Cells.Item(80, 6).Value = "Inactive"

$ws.Cells.Item(81, 1).Value = "2026-01-30"
$ws.Cells.Item(81, 2).Value = "16:55:27"
$ws.Cells.Item(81, 3).Value = "16:00"
$ws.Cells.Item(81, 4).Value = "Bathroom"
$ws.Cells.Item(81, 5).Value = "No Motion"
$ws.Cells.Item(81, 6).Value = "Inactive"

$ws.Cells.Item(82, 1).Value = "2026-01-30"
$ws.Cells.Item(82, 2).Value = "16:55:32"
$ws.Cells.Item(82, 3).Value = "16:00"
$ws.Cells.Item(82, 4).Value = "Bathroom"
$ws.Cells.Item(82, 5).Value = "No Motion"
$ws.Cells.Item(82, 6).Value = "Inactive"

$ws.Cells.Item(83, 1).Value = "2026-01-30"
$ws.Cells.Item(83, 2).Value = "16:55:37"
$ws.Cells.Item(83, 3).Value = "16:00"
$ws.Cells.Item(83, 4).Value = "Bathroom"
$ws.Cells.Item(83, 5).Value = "No Motion"
$ws.Cells.Item(83, 6).Value = "Inactive"

$ws.Cells.Item(84, 1).Value = "2026-01-30"
$ws.Cells.Item(84, 2).Value = "16:55:42"
$ws.Cells.Item(84, 3).Value = "16:00"
$ws.Cells.Item(84, 4).Value = "Bathroom"
$ws.Cells.Item(84, 5).Value = "No Motion"
$ws.Cells.Item(84, 6).Value = "Inactive"

$ws.Cells.Item(85, 1).Value = "2026-01-30"
$ws.Cells.Item(85, 2).Value = "16:55:48"
$ws.Cells.Item(85, 3).Value = "16:00"
$ws.Cells.Item(85, 4).Value = "Bathroom"
$ws.Cells.Item(85, 5).Value = "No Motion"
$ws.Cells.Item(85, 6).Value = "Inactive"

$ws.Cells.Item(86, 1).Value = "2026-01-30"
$ws.Cells.Item(86, 2).Value = "16:55:53"
$ws.Cells.Item(86, 3).Value = "16:00"
$ws.Cells.Item(86, 4).Value = "Bathroom"
$ws.Cells.Item(86, 5).Value = "No Motion"
$ws.Cells.Item(86, 6).Value = "Inactive"

$ws.Cells.Item(87, 1).Value = "2026-01-30"
$ws.Cells.Item(87, 2).Value = "16:55:57"
$ws.Cells.Item(87, 3).Value = "16:00"
$ws.Cells.Item(87, 4).Value = "Bathroom"
$ws.Cells.Item(87, 5).Value = "No Motion"
$ws.Cells.Item(87, 6).Value = "Inactive"

$ws = $wb.Worksheets.Item("Humidity")
$ws.Range("A45:A51").NumberFormat = "@"
$ws.Range("E45:E51").NumberFormat = "@"

$ws.Cells.Item(45, 1).Value = "2026-01-30"
$ws.Cells.Item(45, 2).Value = "16:55:03"
$ws.Cells.Item(45, 3).Value = "16:00"
$ws.Cells.Item(45, 4).Value = "Bathroom"
$ws.Cells.Item(45, 5).Value = "87.8%"
$ws.Cells.Item(45, 6).Value = "Active"

$ws.Cells.Item(46, 1).Value = "2026-01-30"
$ws.Cells.Item(46, 2).Value = "16:55:08"
$ws.Cells.Item(46, 3).Value = "16:00"
$ws.Cells.Item(46, 4).Value = "Bathroom"
$ws.Cells.Item(46, 5).Value = "87.7%"
$ws.Cells.Item(46, 6).Value = "Active"

$ws.Cells.Item(47, 1).Value = "2026-01-30"
$ws.Cells.Item(47, 2).Value = "16:55:23"
$ws.Cells.Item(47, 3).Value = "16:00"
$ws.Cells.Item(47, 4).Value = "Bathroom"
$ws.Cells.Item(47, 5).Value = "87.7%"
$ws.Cells.Item(47, 6).Value = "Active"

$ws.Cells.Item(48, 1).Value = "2026-01-30"
$ws.Cells.Item(48, 2).Value = "16:55:28"
$ws.Cells.Item(48, 3).Value = "16:00"
$ws.Cells.Item(48, 4).Value = "Bathroom"
$ws.Cells.Item(48, 5).Value = "86.9%"
$ws.Cells.Item(48, 6).Value = "Active"

$ws.Cells.Item(49, 1).Value = "2026-01-30"
$ws.Cells.Item(49, 2).Value = "16:55:33"
$ws.Cells.Item(49, 3).Value = "16:00"
$ws.Cells.Item(49, 4).Value = "Bathroom"
$ws.Cells.Item(49, 5).Value = "87.7%"
$ws.Cells.Item(49, 6).Value = "Active"

$ws.Cells.Item(50, 1).Value = "2026-01-30"
$ws.Cells.Item(50, 2).Value = "16:55:43"
$ws.Cells.Item(50, 3).Value = "16:00"
$ws.Cells.Item(50, 4).Value = "Bathroom"
$ws.Cells.Item(50, 5).Value = "87.7%"
$ws.Cells.Item(50, 6).Value = "Active"

$ws.Cells.Item(51, 1).Value = "2026-01-30"
$ws.Cells.Item(51, 2).Value = "16:55:48"
$ws.Cells.Item(51, 3).Value = "16:00"
$ws.Cells.Item(51, 4).Value = "Bathroom"
$ws.Cells.Item(51, 5).Value = "87.7%"
$ws.Cells.Item(51, 6).Value = "Active"

$ws = $wb.Worksheets.Item("Temperature")
$ws.Range("A7:A13").NumberFormat = "@"

$ws.Cells.Item(7, 1).Value = "2026-01-30"
$ws.Cells.Item(7, 2).Value = "16:55:04"
$ws.Cells.Item(7, 3).Value = "16:00"
$ws.Cells.Item(7, 4).Value = "Bathroom"
$ws.Cells.Item(7, 5).Value = "22.6C"
$ws.Cells.Item(7, 6).Value = "Active"

$ws.Cells.Item(8, 1).Value = "2026-01-30"
$ws.Cells.Item(8, 2).Value = "16:55:08"
$ws.Cells.Item(8, 3).Value = "16:00"
$ws.Cells.Item(8, 4).Value = "Bathroom"
$ws.Cells.Item(8, 5).Value = "22.5C"
$ws.Cells.Item(8, 6).Value = "Active"

$ws.Cells.Item(9, 1).Value = "2026-01-30"
$ws.Cells.Item(9, 2).Value = "16:55:23"
$ws.Cells.Item(9, 3).Value = "16:00"
$ws.Cells.Item(9, 4).Value = "Bathroom"
$ws.Cells.Item(9, 5).Value = "22.6C"
$ws.Cells.Item(9, 6).Value = "Active"

$ws.Cells.Item(10, 1).Value = "2026-01-30"
$ws.Cells.Item(10, 2).Value = "16:55:28"
$ws.Cells.Item(10, 3).Value = "16:00"
$ws.Cells.Item(10, 4).Value = "Bathroom"
$ws.Cells.Item(10, 5).Value = "22.7C"
$ws.Cells.Item(10, 6).Value = "Active"

$ws.Cells.Item(11, 1).Value = "2026-01-30"
$ws.Cells.Item(11, 2).Value = "16:55:33"
$ws.Cells.Item(11, 3).Value = "16:00"
$ws.Cells.Item(11, 4).Value = "Bathroom"
$ws.Cells.Item(11, 5).Value = "22.6C"
$ws.Cells.Item(11, 6).Value = "Active"

$ws.Cells.Item(12, 1).Value = "2026-01-30"
$ws.Cells.Item(12, 2).Value = "16:55:43"
$ws.Cells.Item(12, 3).Value = "16:00"
$ws.Cells.Item(12, 4).Value = "Bathroom"
$ws.Cells.Item(12, 5).Value = "22.6C"
$ws.Cells.Item(12, 6).Value = "Active"

$ws.Cells.Item(13, 1).Value = "2026-01-30"
$ws.Cells.Item(13, 2).Value = "16:55:48"
$ws.Cells.Item(13, 3).Value = "16:00"
$ws.Cells.Item(13, 4).Value = "Bathroom"
$ws.Cells.Item(13, 5).Value = "22.6C"
$ws.Cells.Item(13, 6).Value = "Active"

$ws = $wb.Worksheets.Item("Proximity")
$ws.Range("A27:A31").NumberFormat = "@"

$ws.Cells.Item(27, 1).Value = "2026-01-30"
$ws.Cells.Item(27, 2).Value = "16:55:03"
$ws.Cells.Item(27, 3).Value = "16:00"
$ws.Cells.Item(27, 4).Value = "Living Room Main Door"
$ws.Cells.Item(27, 5).Value = "Detected"
$ws.Cells.Item(27, 6).Value = "Active"

$ws.Cells.Item(28, 1).Value = "2026-01-30"
$ws.Cells.Item(28, 2).Value = "16:55:03"
$ws.Cells.Item(28, 3).Value = "16:00"
$ws.Cells.Item(28, 4).Value = "Living Room Main Door"
$ws.Cells.Item(28, 5).Value = "Clear"
$ws.Cells.Item(28, 6).Value = "Inactive"

$ws.Cells.Item(29, 1).Value = "2026-01-30"
$ws.Cells.Item(29, 2).Value = "16:55:11"
$ws.Cells.Item(29, 3).Value = "16:00"
$ws.Cells.Item(29, 4).Value = "Living Room Main Door"
$ws.Cells.Item(29, 5).Value = "Detected"
$ws.Cells.Item(29, 6).Value = "Active"

$ws.Cells.Item(30, 1).Value = "2026-01-30"
$ws.Cells.Item(30, 2).Value = "16:55:13"
$ws.Cells.Item(30, 3).Value = "16:00"
$ws.Cells.Item(30, 4).Value = "Living Room Main Door"
$ws.Cells.Item(30, 5).Value = "Clear"
$ws.Cells.Item(30, 6).Value = "Inactive"

$ws.Cells.Item(31, 1).Value = "2026-01-30"
$ws.Cells.Item(31, 2).Value = "16:55:16"
$ws.Cells.Item(31, 3).Value = "16:00"
$ws.Cells.Item(31, 4).Value = "Living Room Main Door"
$ws.Cells.Item(31, 5).Value = "Detected"
$ws.Cells.Item(31, 6).Value = "Active"

$ws = $wb.Worksheets.Item("mmWave")
$ws.Range("A16:A21").NumberFormat = "@"

$ws.Cells.Item(16, 1).Value = "2026-01-30"
$ws.Cells.Item(16, 2).Value = "16:55:03"
$ws.Cells.Item(16, 3).Value = "16:00"
$ws.Cells.Item(16, 4).Value = "Living Room"
$ws.Cells.Item(16, 5).Value = "PRESENCE_DETECTED"
$ws.Cells.Item(16, 6).Value = "Active"

$ws.Cells.Item(17, 1).Value = "2026-01-30"
$ws.Cells.Item(17, 2).Value = "16:55:13"
$ws.Cells.Item(17, 3).Value = "16:00"
$ws.Cells.Item(17, 4).Value = "Living Room"
$ws.Cells.Item(17, 5).Value = "PRESENCE_DETECTED"
$ws.Cells.Item(17, 6).Value = "Active"

$ws.Cells.Item(18, 1).Value = "2026-01-30"
$ws.Cells.Item(18, 2).Value = "16:55:23"
$ws.Cells.Item(18, 3).Value = "16:00"
$ws.Cells.Item(18, 4).Value = "Living Room"
$ws.Cells.Item(18, 5).Value = "PRESENCE_DETECTED"
$ws.Cells.Item(18, 6).Value = "Active"

$ws.Cells.Item(19, 1).Value = "2026-01-30"
$ws.Cells.Item(19, 2).Value = "16:55:34"
$ws.Cells.Item(19, 3).Value = "16:00"
$ws.Cells.Item(19, 4).Value = "Living Room"
$ws.Cells.Item(19, 5).Value = "PRESENCE_DETECTED"
$ws.Cells.Item(19, 6).Value = "Active"

$ws.Cells.Item(20, 1).Value = "2026-01-30"
$ws.Cells.Item(20, 2).Value = "16:55:48"
$ws.Cells.Item(20, 3).Value = "16:00"
$ws.Cells.Item(20, 4).Value = "Living Room"
$ws.Cells.Item(20, 5).Value = "PRESENCE_DETECTED"
$ws.Cells.Item(20, 6).Value = "Active"

$ws.Cells.Item(21, 1).Value = "2026-01-30"
$ws.Cells.Item(21, 2).Value = "16:55:58"
$ws.Cells.Item(21, 3).Value = "16:00"
$ws.Cells.Item(21, 4).Value = "Living Room"
$ws.Cells.Item(21, 5).Value = "PRESENCE_DETECTED"
$ws.Cells.Item(21, 6).Value = "Active"
